$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 56.69
$ws.Range("B4").Value = 42.06
$ws.Range("B5").Value = 62.07
$ws.Range("B7").Value = 49.36
$ws.Range("B8").Value = 67.66
$ws.Range("B10").Value = 57.26
$ws.Range("B12").Value = 43.56
$ws.Range("B14").Value = 58.26
$ws.Range("B16").Value = 18.49
$ws.Range("B17").Value = 56.86
$ws.Range("B19").Value = 45.59
$ws.Range("B20").Value = 63.76
